# Update the worker / period / value data in the account-statement table
# (rows 16-20) to reflect the reordered / corrected records described in
# the commit: previous statements removed and new ones added, database
# updated.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 16: KETTY LUZ ACOSTA MARTINEZ / period 2304
$ws.Range("C16").Value = "1049452250"
$ws.Range("D16").Value = "KETTY LUZ ACOSTA MARTINEZ"
$ws.Range("E16").Value = "2304"
$ws.Range("F16").Value = 46400

# Row 17: KETTY LUZ ACOSTA MARTINEZ / period 2303
$ws.Range("C17").Value = "1049452250"
$ws.Range("D17").Value = "KETTY LUZ ACOSTA MARTINEZ"
$ws.Range("E17").Value = "2303"
$ws.Range("F17").Value = 46400

# Row 18: EDELMIRA CASTILLO ALVAREZ / period 2303
$ws.Range("C18").Value = "45372092"
$ws.Range("D18").Value = "EDELMIRA CASTILLO ALVAREZ"
$ws.Range("E18").Value = "2303"
$ws.Range("F18").Value = 13920

# Row 19: NORELIS MENDOZA ROCHA / period 2308
$ws.Range("C19").Value = "1049941768"
$ws.Range("D19").Value = "NORELIS MENDOZA ROCHA"
$ws.Range("E19").Value = "2308"
$ws.Range("F19").Value = 6187

# Row 20: NORELIS MENDOZA ROCHA / period 2307
$ws.Range("C20").Value = "1049941768"
$ws.Range("D20").Value = "NORELIS MENDOZA ROCHA"
$ws.Range("E20").Value = "2307"
$ws.Range("F20").Value = 46400
